# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handback DateTime"
# timestamps on the Overview, zh-cn and de-de sheets to reflect a new report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first data row
$wsOverview.Range("G2").Value = "2016-08-30 13:10:14"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for the first data row
$wsZhCn.Range("H2").Value = "2016-08-30 13:09:58"
$wsZhCn.Range("K2").Value = "2016-08-30 13:10:37"

# de-de sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for the first data row
$wsDeDe.Range("H2").Value = "2016-08-30 13:10:14"
$wsDeDe.Range("K2").Value = "2016-08-30 13:10:44"
